$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing the choke resistor (Mfg Part # ERJ-1GNF3000C, Ref R8)
# and delete the entire row, shifting all rows below it up by one.
$target = $ws.Range("B2:B32").Find("ERJ-1GNF3000C")
if ($target -ne $null) {
    $row = $target.Row
    $ws.Rows.Item($row).Delete()
}

# Update the active selection to match the post-edit state
$ws.Range("D20").Select()
